$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-12 Saturday", "2024-10-13 Sunday"),
    @("994÷2=", "898÷7="),
    @("242÷5=", "674÷5="),
    @("645÷8=", "377÷5="),
    @("767÷8=", "412÷2="),
    @("679÷9=", "920÷8="),
    @("555÷7=", "332÷9="),
    @("432÷5=", "360÷8="),
    @("611÷2=", "277÷8="),
    @("240÷8=", "376÷2="),
    @("956÷7=", "198÷9="),
    @("237÷9=", "887÷5="),
    @("356÷8=", "160÷6="),
    @("298÷4=", "562÷3="),
    @("781÷8=", "326÷6="),
    @("123÷4=", "127÷7="),
    @("267÷8=", "570÷6="),
    @("926÷8=", "651÷4="),
    @("365÷3=", "100÷2="),
    @("150÷2=", "289÷2="),
    @("255÷5=", "649÷6="),
    @("559÷9=", "817÷9="),
    @("248÷5=", "145÷2="),
    @("884÷6=", "652÷4="),
    @("730÷3=", "363÷7="),
    @("558÷6=", "652÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
